# Update the cached "datetimeFigureOut" date field text from 10/26/2024 to
# 5/30/2025 everywhere it appears: the slide master, every slide layout
# (CustomLayout) hanging off the master, and the notes master.
#
# PowerPoint normally re-caches these automatic date fields on save; here we
# walk every placeholder that currently shows the stale date and set its text
# to the new one.

$p = $ppt.ActivePresentation

$oldDate = "10/26/2024"
$newDate = "5/30/2025"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1) Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# 2) Every slide layout under the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# 3) Notes master
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes
